$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date in column C for rows 2-6 from 45208 to 45212
$ws.Range("C2:C6").Value = 45212

# Update hyperlink formulas in row 2 (S2:Y2) to include the updated file names
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/artfynd/A 32298-2023 artfynd.xlsx", "A 32298-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/kartor/A 32298-2023 karta.png", "A 32298-2023")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/knärot/A 32298-2023 karta knärot.png", "A 32298-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/klagomål/A 32298-2023 fsc-klagomål.docx", "A 32298-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/klagomålsmail/A 32298-2023 fsc-klagomål mail.docx", "A 32298-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/tillsyn/A 32298-2023 tillsynsbegäran.docx", "A 32298-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/ti,llsynsmail/A 32298-2023 tillsynsbegäran mail.docx", "A 32298-2023")'
